# Update countries & provincias Spain
# - Move "Bolivia" earlier in the country list (now appears right before
#   "Estado de Palestina" instead of right before "Sri Lanka"), and refresh
#   the statistics block for the affected rows accordingly.
# - Update the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country name + stats for rows 106-112 (Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$rows = @(
    @{ Row = 106; Pais = "Bolivia";             B = 264; C = 54; D = 2;   E = 244; F = 3; G = 3; H = 18 },
    @{ Row = 107; Pais = "Estado de Palestina";  B = 263; C = 0;  D = 44;  E = 218; F = 0; G = 0; H = 1  },
    @{ Row = 108; Pais = "Vietnam";              B = 251; C = 0;  D = 126; E = 125; F = 8; G = 0; H = 0  },
    @{ Row = 109; Pais = "Montenegro";           B = 248; C = 0;  D = 4;   E = 242; F = 7; G = 0; H = 2  },
    @{ Row = 110; Pais = "Senegal";              B = 244; C = 0;  D = 113; E = 129; F = 1; G = 0; H = 2  },
    @{ Row = 111; Pais = "Banglades";            B = 218; C = 0;  D = 33;  E = 165; F = 1; G = 0; H = 20 },
    @{ Row = 112; Pais = "Georgia";              B = 211; C = 0;  D = 50;  E = 158; F = 6; G = 0; H = 3  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Pais
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}

# Update timestamp footer text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 06:22"
